$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "percent err" formula in D2:D7 to be *100 (a true percentage) ---
$ws.Range("D2").Formula = "=100*(C2 - B2)/B2"
$ws.Range("D3").Formula = "=100*(C3 - B3)/B3"
$ws.Range("D4:D7").Formula = "=100*(C4 - B4)/B4"

# B5 is re-entered individually (breaks away from the B3:B7 shared formula group)
$ws.Range("B5").Formula = "=A5/3.141596"

# Corrected measured diameter for the last sample
$ws.Range("C7").Value = 28.29

# --- New section: DBH estimate using segmentation ("w/ the whole process for each sample") ---
$ws.Range("A11").Value = "w/ the whole process for each sample"

$ws.Range("B12").Value = "measured"
$ws.Range("C12").Value = "computed"
$ws.Range("D12").Value = "percent err"

# Row 13 (mirrors row 2)
$ws.Range("A13").Value = 140
$ws.Range("B13").Formula = "=A13/3.141596"
$ws.Range("C13").Value = 0
$ws.Range("D13").Formula = "=100*(C13 - B13)/B13"

# Row 14 (mirrors row 3)
$ws.Range("A14").Value = 67
$ws.Range("C14").Value = 19.119
$ws.Range("C14").Font.Name = "Arial Unicode MS"
$ws.Range("C14").Font.Size = 10
$ws.Range("D14").Formula = "=100*(C14 - B14)/B14"

# Row 15 (mirrors row 4)
$ws.Range("A15").Value = 50.5
$ws.Range("C15").Value = 15.936
$ws.Range("C15").Font.Name = "Arial Unicode MS"
$ws.Range("C15").Font.Size = 10

# B14:B15 entered together as one shared-formula fill
$ws.Range("B14:B15").Formula = "=A14/3.141596"

# Row 16 (mirrors row 5) - re-entered individually, breaking away from the B group
$ws.Range("A16").Value = 12.5
$ws.Range("B16").Formula = "=A16/3.141596"
$ws.Range("C16").Value = 3.87
$ws.Range("C16").Font.Name = "Arial Unicode MS"
$ws.Range("C16").Font.Size = 10

# D15:D18 entered together as one shared-formula fill
$ws.Range("D15:D18").Formula = "=100*(C15 - B15)/B15"

# Row 17 (mirrors row 6)
$ws.Range("A17").Value = 146.3
$ws.Range("C17").Value = 43.494
$ws.Range("C17").Font.Name = "Arial Unicode MS"
$ws.Range("C17").Font.Size = 10

# Row 18 (mirrors row 7, using the corrected measurement)
$ws.Range("A18").Value = 92.8
$ws.Range("C18").Value = 28.29
$ws.Range("C18").Font.Name = "Arial Unicode MS"
$ws.Range("C18").Font.Size = 10

# B17:B18 entered together as one shared-formula fill
$ws.Range("B17:B18").Formula = "=A17/3.141596"

# Move / update the active selection to A11, matching the saved view state
$ws.Range("A11").Select()
